# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The workbook is a "Estado de Cuenta" (account statement) template for
# NIT 9018095605. This edit swaps in a new worker / first overdue period
# ("parte 1") and removes the old worker's stale multi-period rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the single data row (row 16) with the new worker/period ---
$ws.Range("C16").Value = "1047483592"
$ws.Range("D16").Value = "CESAR ANDRES REDONDO MANJARREZ"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# --- Update the summary header fields to match the new single period ---
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1

# Column D ("Nombre Trabajador") is best-fit width; let Excel resize it now
# that the name text is shorter.
$ws.Columns("D:D").AutoFit()

# --- Remove the now-obsolete extra period rows (2506..2502) ---
# These were rows 17-21 in the old layout; deleting them shifts the
# signature block (previously rows 26-27) up to rows 21-22.
$ws.Rows("17:21").Delete()

"Updated worker record, summary totals, and removed stale period rows."
